$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shared-string order in the target workbook is Max(0), Min(1), Avg(2) -
# write the labels in that order so the generated sharedStrings table
# matches, then fill in the formulas.
$ws.Range("E4").Value = "Max"
$ws.Range("E1").Value = "Min"
$ws.Range("E7").Value = "Avg"

# Min block (rows 1-2)
$ws.Range("E2").Formula = "=MIN(A1:A10)/200000000*1000"
$ws.Range("F2").Formula = "=MIN(C1:C10)/200000000*1000"

# Max block (rows 4-5)
$ws.Range("E5").Formula = "=MAX(A1:A10)/200000000*1000"
$ws.Range("F5").Formula = "=MAX(C1:C10)/200000000*1000"

# Avg block (rows 7-8)
$ws.Range("E8").Formula = "=AVERAGE(A1:A10)/200000000*1000"
$ws.Range("F8").Formula = "=AVERAGE(C1:C10)/200000000*1000"

# Remove the old single average formula in row 11
$ws.Range("C11").ClearContents()

# Update selection to match the new active cell
$ws.Range("E9").Select()

# Add page setup (portrait orientation) matching the new sheet1.xml
$ws.PageSetup.Orientation = 1
